# Updated cryptos list on Mon Jan  1 18:20:25 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of column letters (as used in the source data) to column indices.
$colmap = @{ B = 2; C = 3; D = 4; E = 5 }

# Each hashtable below describes the cells that changed for a given row.
# Only the columns present in the hashtable are touched; everything else
# is left exactly as it was.
$updates = @(
    @{Row=2; D='43.357.63'; E='  +1.78%  '},
    @{Row=3; D='2.329.08'; E='  +1.48%  '},
    @{Row=4; E='  +0.00%  '},
    @{Row=5; D='313.22'; E='  -0.83%  '},
    @{Row=6; D='108.48'; E='  +4.22%  '},
    @{Row=7; D='0.630'; E='  +1.03%  '},
    @{Row=8; E='  +0.18%  '},
    @{Row=9; D='0.612'; E='  +1.71%  '},
    @{Row=10; D='40.50'; E='  +2.82%  '},
    @{Row=11; D='0.0918'; E='  +1.27%  '},
    @{Row=12; D='8.50'; E='  +0.11%  '},
    @{Row=13; E='  -1.07%  '},
    @{Row=14; D='1.00'; E='  +0.13%  '},
    @{Row=15; D='15.41'; E='  +0.59%  '},
    @{Row=16; D='2.675.41'; E='  +1.22%  '},
    @{Row=17; D='2.324.88'; E='  +1.07%  '},
    @{Row=18; D='43.248.21'; E='  +1.55%  '},
    @{Row=19; D='7.52'; E='  -0.17%  '},
    @{Row=20; E='  +0.68%  '},
    @{Row=21; E='  -12.31%  '},
    @{Row=22; D='73.91'; E='  -0.23%  '},
    @{Row=23; D='3.52'; E='  -0.70%  '},
    @{Row=24; D='267.97'; E='  +1.77%  '},
    @{Row=25; D='2.26'; E='  +2.31%  '},
    @{Row=26; E='  +0.03%  '},
    @{Row=27; D='7.68'; E='  +11.40%  '},
    @{Row=28; D='11.11'; E='  +2.21%  '},
    @{Row=29; E='  -2.07%  '},
    @{Row=30; D='38.99'; E='  +4.57%  '},
    @{Row=31; D='22.59'; E='  +1.12%  '},
    @{Row=32; D='166.90'; E='  +0.14%  '},
    @{Row=33; D='0.0882'; E='  +0.86%  '},
    @{Row=34; D='2.76'; E='  +6.06%  '},
    @{Row=35; D='0.131'; E='  +0.15%  '},
    @{Row=36; D='4.72'; E='  +3.04%  '},
    @{Row=37; D='0.113'; E='  -1.74%  '},
    @{Row=38; D='0.0360'; E='  +2.77%  '},
    @{Row=39; D='2.84'; E='  +5.75%  '},
    @{Row=40; D='3.73'; E='  +0.87%  '},
    @{Row=41; D='1.65'; E='  +5.12%  '},
    @{Row=42; D='104.46'; E='  +11.44%  '},
    @{Row=43; B='MultiversX'; C='https://coinranking.com/coin/omwkOTglq+multiversx-egld'; D='71.28'; E='  +2.67%  '},
    @{Row=44; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.235'; E='  +2.42%  '},
    @{Row=45; D='13.30'; E='  +7.83%  '},
    @{Row=46; D='1.00'; E='  +0.28%  '},
    @{Row=47; D='114.03'; E='  +0.01%  '},
    @{Row=48; B='ordi'; C='https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'; D='77.77'; E='  -2.96%  '},
    @{Row=49; B='TheGraph'; C='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; D='0.219'; E='  +17.63%  '},
    @{Row=50; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.663.25'; E='  -3.12%  '},
    @{Row=51; B='THORChain'; C='https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; D='5.31'; E='  +5.24%  '}
)

foreach ($item in $updates) {
    $r = $item.Row
    foreach ($key in $item.Keys) {
        if ($key -ne "Row") {
            $col = $colmap[$key]
            $cell = $ws.Cells.Item($r, $col)
            # Force text formatting so numeric-looking strings (e.g. "1.00",
            # "43.357.63", "0.0360") keep their exact original formatting
            # instead of being coerced into Excel numbers.
            $cell.NumberFormat = "@"
            $cell.Value = $item[$key]
        }
    }
}
